$d = $word.ActiveDocument

# 1. Remove the "[Grade]" placeholder run entirely (leaves an empty paragraph)
$d.Content.Find.Execute("[Grade]", $false, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# 2. Collapse the 3-run "[" + "Dept Head" + "]" into a single run reading "[Department Head]"
$d.Content.Find.Execute("[Dept Head]", $false, $false, $false, $false, $false, $true, 1, $false, "[Department Head]", 2) | Out-Null

# 3. Collapse the 3-run "[" + "College Dean" + "]" into a single run reading "[College Dean]"
$d.Content.Find.Execute("[College Dean]", $false, $false, $false, $false, $false, $true, 1, $false, "[College Dean]", 2) | Out-Null

# 4. Turn the plain "Department Head" signature line into a misspelled "Deparment Head",
#    split across two runs and flagged with proofErr spell-check markers, matching the
#    structure Word produces when its proofing engine flags a word it doesn't recognize.
$n = $d.Paragraphs.Count
for ($i = 1; $i -le $n; $i++) {
    $para = $d.Paragraphs($i)
    $coreText = $para.Range.Text.TrimEnd([char]13, [char]7)
    if ($coreText -eq "Department Head") {
        $r = $para.Range
        $xml = "<?xml version='1.0' encoding='UTF-8' standalone='yes'?><pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' xmlns:w14='http://schemas.microsoft.com/office/word/2010/wordml'><w:body><w:p w14:paraId='1F344026' w14:textId='59F22672' w:rsidR='00897F0B' w:rsidRPr='00DA048A' w:rsidRDefault='00897F0B' w:rsidP='00897F0B'><w:pPr><w:spacing w:before='100' w:beforeAutospacing='1' w:after='100' w:afterAutospacing='1'/><w:jc w:val='center'/><w:rPr><w:rFonts w:ascii='Tahoma' w:eastAsia='Times New Roman' w:hAnsi='Tahoma' w:cs='Tahoma'/><w:color w:val='191919'/><w:sz w:val='20'/><w:szCs w:val='20'/></w:rPr></w:pPr><w:proofErr w:type='spellStart'/><w:r><w:rPr><w:rFonts w:ascii='Tahoma' w:eastAsia='Times New Roman' w:hAnsi='Tahoma' w:cs='Tahoma'/><w:color w:val='191919'/><w:sz w:val='20'/><w:szCs w:val='20'/></w:rPr><w:t>Deparment</w:t></w:r><w:proofErr w:type='spellEnd'/><w:r><w:rPr><w:rFonts w:ascii='Tahoma' w:eastAsia='Times New Roman' w:hAnsi='Tahoma' w:cs='Tahoma'/><w:color w:val='191919'/><w:sz w:val='20'/><w:szCs w:val='20'/></w:rPr><w:t xml:space='preserve'> Head</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
        $r.InsertXML($xml) | Out-Null
        break
    }
}
